# Trans-Tasman charts updated to 29 November
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert 7 new rows above row 2 (pushes existing data down by 7 rows) ---
$insertRange = $ws.Range("A2:B8")
$insertRange.EntireRow.Insert()

# --- New data for the 7 inserted rows (newest first) ---
$newDates  = @(44529, 44528, 44527, 44526, 44525, 44524, 44523)
$newValues = @(6366060, 6365073, 6362348, 6355156, 6348980, 6341726, 6334061)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newValues[$i]
    $ws.Rows.Item($r).RowHeight = 18
    $ws.Cells.Item($r, 1).NumberFormat = "d-mmm-yy"
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
}

# --- Fill in values that were previously blank, now shifted to rows 237-243 ---
$fillValues = @(38445, 35475, 31795, 31212, 31212, 30383, 29796)
for ($i = 0; $i -lt $fillValues.Length; $i++) {
    $r = 237 + $i
    $ws.Cells.Item($r, 2).Value = $fillValues[$i]
}

# --- Update selection ---
$ws.Range("D9").Select()

# --- Make the "second doses" header font red ---
$ws.Range("B1").Font.Color = 255
